# Fill in the remaining patient rows on the tcia-submission-template sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - TCGA patient
$ws.Range("A4").Value = 7980790
$ws.Range("B4").Value = "TCGA-12-5436"
$ws.Range("C4").Value = "TCGA-BLCA"
$ws.Range("D4").Value = "4/9/2007"

# Row 5 - CPTAC patient
$ws.Range("A5").Value = 12348975
$ws.Range("B5").Value = "CPTAC-UCEC-0001"
$ws.Range("C5").Value = "CPTAC-Uterine"
$ws.Range("D5").Value = "3/31/2017"

# Row 6 - CPTAC patient
$ws.Range("A6").Value = 13497812
$ws.Range("B6").Value = "CPTAC-UCEC-0002"
$ws.Range("C6").Value = "CPTAC-Uterine"
$ws.Range("D6").Value = "3/31/2017"

# Move the active selection to D6, matching where the user ended up editing.
$ws.Range("D6").Select()
